# Apply the edit: insert a new data row at row 42 (pushing existing rows
# 42-124 down to 43-125) and populate the new row 42 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 42; this shifts rows 42:124
# down to 43:125 and leaves row 42 empty for the new record.
$ws.Rows("42:42").Insert()

# Populate the new row 42 with its data.
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C42").Value = "Arica y Parinacota"
$ws.Range("D42").Value = 44965
$ws.Range("D42").NumberFormat = $ws.Range("D43").NumberFormat
$ws.Range("E42").Value = 15
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100102
$ws.Range("H42").Value = "Cítricos"
$ws.Range("I42").Value = 100102005
$ws.Range("J42").Value = "Naranja"
$ws.Range("K42").Value = "Lane Late"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 250
$ws.Range("N42").Value = 900
$ws.Range("O42").Value = 950
$ws.Range("P42").Value = 920
$ws.Range("Q42").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R42").Value = "Región de Coquimbo"
$ws.Range("S42").Value = 920
$ws.Range("T42").Value = 1
